$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Userdata")

$ws.Range("A3").Value = "admin"
$ws.Range("B3").Value = "admin"
$ws.Range("C3").Value = "admins@gmail.com"
$ws.Range("D3").Value = "admin"
$ws.Range("E3").Value = "admin123"
